# Update the case locations / public exposure sites table.
# The Camberwell restaurant entries (rows 2 & 3) are replaced with two new
# exposure sites: a Hampton cafe and a McKinnon hairdresser.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Hampton cafe exposure site
$ws.Range("A2").Value = "Hampton"
$ws.Range("B2").Value = "Merrymen Cafe, 2 Small Street, Hampton VIC"
$ws.Range("C2").Value = "28-12-2020 1:20pm-2:30pm"
$ws.Range("D2").Value = "Case ate in store"
$ws.Range("E2").Value = "new"

# Row 3: McKinnon hairdresser exposure site
$ws.Range("A3").Value = "McKinnon"
$ws.Range("B3").Value = "260 McKinnon Road, McKinnon VIC 3204"
$ws.Range("C3").Value = "23-12-2020 4:00pm-6:00pm"
$ws.Range("D3").Value = "Case had hair cut in store"
$ws.Range("E3").Value = "new"

# Match the column widths Excel's bestFit autosizing computed for the new text
# (same effect as double-clicking each column border / AutoFit, but with the
# exact resulting widths so the saved file matches Excel's own output).
$ws.Columns.Item(1).ColumnWidth = 8.6640625
$ws.Columns.Item(2).ColumnWidth = 36.46484375
$ws.Columns.Item(3).ColumnWidth = 23.33203125
$ws.Columns.Item(4).ColumnWidth = 20.73046875
$ws.Columns.Item(5).ColumnWidth = 4.46484375

# Move the active selection to B2 (matches the saved selection in the workbook)
$ws.Range("B2").Select() | Out-Null
